$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values differ between row 3 and row 4 and must be swapped
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $r3 = $ws.Range($col + "3")
    $r4 = $ws.Range($col + "4")
    $v3 = $r3.Value()
    $v4 = $r4.Value()
    $r3.Value = $v4
    $r4.Value = $v3
}

# AC3/AC4: the "På sälg" comment moves from row 4 to row 3
$ws.Range("AC3").Value = "På sälg"
$ws.Range("AC4").Value = $null
